# Update header text in cell A1 from "Company" to "company"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "company"

# Move the active selection from E2 to A2
$ws.Range("A2").Select()
